# Update Section_A and Section_B timetable sheets with the new schedule values.
$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("Section_A")
$wsB = $wb.Worksheets.Item("Section_B")

# --- Section_A (sheet1) ---
$wsA.Range("B2").Value = "Free"
$wsA.Range("C2").Value = "Free"
$wsA.Range("D2").Value = "Free"
$wsA.Range("E2").Value = "Free"
$wsA.Range("F2").Value = "CS261"

$wsA.Range("B3").Value = "CS263"
$wsA.Range("C3").Value = "CS261"
$wsA.Range("D3").Value = "CS264"
$wsA.Range("E3").Value = "Free"
$wsA.Range("F3").Value = "Free"

$wsA.Range("B5").Value = "MA261"
$wsA.Range("C5").Value = "CS264"
$wsA.Range("D5").Value = "Free"
$wsA.Range("E5").Value = "Free"
$wsA.Range("F5").Value = "CS264"

$wsA.Range("B6").Value = "CS261"
$wsA.Range("C6").Value = "CS263"
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "Free"
$wsA.Range("F6").Value = "Free"

$wsA.Range("B7").Value = "Free"
$wsA.Range("C7").Value = "Free"
$wsA.Range("D7").Value = "MA261"
$wsA.Range("E7").Value = "Free"
$wsA.Range("F7").Value = "CS263"

# --- Section_B (sheet2) ---
$wsB.Range("B2").Value = "CS264"
$wsB.Range("C2").Value = "Free"
$wsB.Range("D2").Value = "CS263"
$wsB.Range("E2").Value = "CS261"
$wsB.Range("F2").Value = "MA261"

$wsB.Range("B3").Value = "CS263"
$wsB.Range("C3").Value = "Free"
$wsB.Range("D3").Value = "Free"
$wsB.Range("E3").Value = "Free"
$wsB.Range("F3").Value = "CS261"

$wsB.Range("B5").Value = "MA261"
$wsB.Range("C5").Value = "Free"
$wsB.Range("D5").Value = "Free"
$wsB.Range("E5").Value = "Free"
$wsB.Range("F5").Value = "Free"

$wsB.Range("B6").Value = "Free"
$wsB.Range("C6").Value = "CS263"
$wsB.Range("D6").Value = "Free"
$wsB.Range("E6").Value = "CS264"
$wsB.Range("F6").Value = "CS264"

$wsB.Range("B7").Value = "CS261"
$wsB.Range("C7").Value = "Free"
$wsB.Range("D7").Value = "Free"
$wsB.Range("E7").Value = "Free"
$wsB.Range("F7").Value = "Free"
